$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---------------------------------------------------------
# A1, B1 keep their existing text ("Test Case ID" / "Test Case Description").
# C1 used to read "Step Name"; it now shows "Expected Result " (the text
# that used to live in F1). D1 keeps its bold style but loses its text
# ("Test Data" -> blank). The old E1:G1 headers ("Step Description",
# "Expected Result ", "Actual") are removed entirely together with their
# columns.
$ws.Range("C1").Value = "Expected Result "
$ws.Range("D1").ClearContents()
$ws.Range("E1:G1").EntireColumn.Delete()

# --- New test-case rows --------------------------------------------------
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "New customer without new card "
$ws.Range("C2").Value = "No discount"

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "New customer with new card"
$ws.Range("C3").Value = "15% discount"

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Existing customer with no card"
$ws.Range("C4").Value = "No discount"

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Existing customer with card"
$ws.Range("C5").Value = "10% discount"

$ws.Range("A6").Value = 5
$ws.Range("C6").Value = "20% discount"
$ws.Range("B6").Value = "If the new customer has discount coupon and no card"

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "If the new customer has discount coupon and card"
$ws.Range("C7").NumberFormat = "0%"
$ws.Range("C7").Value = "30% discound"

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "If the new customer has no discount coupon and has card"
$ws.Range("C8").Value = "10% discount"

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "If the new customer has no discount coupon and has no card"
$ws.Range("C9").Value = "No discount"

# --- Row height for row 2 -------------------------------------------------
$ws.Rows.Item(2).RowHeight = 20

# --- Column widths ---------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 51.67
$ws.Columns.Item(3).ColumnWidth = 24.5
$ws.Columns.Item(4).ColumnWidth = 23.83

# --- Selection -------------------------------------------------------------
$ws.Range("C10").Select() | Out-Null
